$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Excel auto-converts plain numeric-looking strings (e.g. "1.00", "284.90")
    # to actual numbers, which would silently drop significant trailing zeros.
    # Forcing a Text number format before the write keeps the literal string,
    # then resetting the style back to Normal avoids leaving a stray cell format
    # behind once the value itself is safely stored as text.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '60.648.34'
$ws.Range('E2').Value = '  +2.69%  '
$ws.Range('D3').Value = '2.701.96'
$ws.Range('E3').Value = '  +2.69%  '
Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  +0.11%  '
Set-TextValue $ws.Range('D5') '526.47'
$ws.Range('E5').Value = '  +1.50%  '
Set-TextValue $ws.Range('D6') '144.83'
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  +2.32%  '
$ws.Range('D9').Value = '2.727.13'
$ws.Range('E9').Value = '  +2.66%  '
Set-TextValue $ws.Range('D10') '6.69'
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('E13').Value = '  +2.97%  '
$ws.Range('D14').Value = '3.180.52'
$ws.Range('E14').Value = '  +2.22%  '
$ws.Range('D15').Value = '60.657.44'
$ws.Range('E15').Value = '  +2.78%  '
$ws.Range('D16').Value = '2.898.51'
$ws.Range('E16').Value = '  +8.92%  '
Set-TextValue $ws.Range('D17') '21.32'
$ws.Range('E17').Value = '  +1.58%  '
Set-TextValue $ws.Range('D19') '348.02'
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('E20').Value = '  -0.09%  '
Set-TextValue $ws.Range('D21') '10.62'
$ws.Range('E21').Value = '  +2.75%  '
Set-TextValue $ws.Range('D22') '6.41'
$ws.Range('E22').Value = '  +3.87%  '
Set-TextValue $ws.Range('D23') '0.998'
$ws.Range('E23').Value = '  +0.03%  '
Set-TextValue $ws.Range('D24') '63.64'
$ws.Range('E24').Value = '  +2.54%  '
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('E26').Value = '  +4.63%  '
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('E28').Value = '  +1.54%  '
Set-TextValue $ws.Range('D29') '7.31'
$ws.Range('E29').Value = '  +2.54%  '
Set-TextValue $ws.Range('D30') '6.77'
$ws.Range('E30').Value = '  +8.28%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  +1.63%  '
Set-TextValue $ws.Range('D33') '19.14'
$ws.Range('E33').Value = '  +0.73%  '
Set-TextValue $ws.Range('D34') '150.33'
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('E35').Value = '  +5.37%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D36') '1.23'
$ws.Range('E36').Value = '  +7.72%  '
$ws.Range('B37').Value = 'SuiNetwork'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range('D37') '0.943'
$ws.Range('E37').Value = '  -2.05%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D38') '1.53'
$ws.Range('E38').Value = '  +7.87%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D39') '0.877'
$ws.Range('E39').Value = '  +3.98%  '
Set-TextValue $ws.Range('D40') '37.06'
$ws.Range('E40').Value = '  +0.89%  '
$ws.Range('E41').Value = '  -0.85%  '
Set-TextValue $ws.Range('D42') '284.90'
$ws.Range('E42').Value = '  +2.91%  '
Set-TextValue $ws.Range('D43') '20.17'
$ws.Range('E43').Value = '  +2.78%  '
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('E45').Value = '  +0.64%  '
Set-TextValue $ws.Range('D46') '0.996'
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('D47').Value = '2.141.06'
$ws.Range('E47').Value = '  +7.49%  '
$ws.Range('E48').Value = '  +3.74%  '
$ws.Range('E49').Value = '  +2.36%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range('D50') '10.48'
$ws.Range('E50').Value = '  +1.70%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D51') '4.80'
$ws.Range('E51').Value = '  +3.27%  '
